$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Testmail #5: Wil je deze klant bellen?"
$ws.Range("B6").Value = "Geachte heer/mevrouw,`nBedankt voor uw e-mail. We zullen de klant zo snel mogelijk contacteren. Mocht u nog meer informatie hebben die u met ons wilt delen, dan horen we dat graag.`nMet vriendelijke groet,`n[Naam bedrijf]"
$ws.Range("C6").Value = "Wil je deze klant bellen?"
$ws.Range("D6").Value = "mailmind.test@zohomail.eu"
$ws.Range("E6").Value = "Klantenservice / Contact"
$ws.Range("F6").Value = "2025-08-04 20:36:10"
$ws.Range("G6").Value = "Ja"
$ws.Range("H6").Value = "Nee"
$ws.Range("I6").Value = "Ja"
$ws.Range("J6").Value = "Nee"

# Remove the automatic row-height expansion caused by the multi-line
# value in B6 so the row keeps using the sheet's default row height.
$ws.Rows.Item(6).AutoFit()
